# Auto-generated script to update cryptos worksheet values
# Applies the per-cell text updates described by the diff, preserving
# each cell as a text value (matching the original inlineStr typing).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "61.107.77"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.650.77"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.61"
$ws.Range("E5").Value = "  +4.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.00"
$ws.Range("E6").Value = "  +0.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E8").Value = "  +0.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.64"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.352"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "3.111.45"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").Value = "61.109.01"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.11"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("E16").Value = "  +2.57%  "
$ws.Range("D17").Value = "2.642.33"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "355.97"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.70"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.26"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.71"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.432"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "0.0₃0861"
$ws.Range("E27").Value = "  +2.12%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.20"
$ws.Range("E30").Value = "  +7.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.57"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.23"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  +3.52%  "
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.924"
$ws.Range("E36").Value = "  +9.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.903"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "310.14"
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("E40").Value = "  +1.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.647"
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0564"
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.06"
$ws.Range("E45").Value = "  +3.86%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.93"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.24"
$ws.Range("E48").Value = "  +7.79%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "1.989.96"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("E51").Value = "  +2.72%  "
